# Append the 2023-06-19 09:30 resale-numbers update as new row 57.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# Text-like columns (Date, Time, Weekday, Week) must stay as plain text,
# matching the rest of the sheet - otherwise Excel auto-converts values
# like "2023-06-19" / "09:29:46" / "25" into a date/time/number.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2023-06-19"
$ws.Range("B$row").Value = "09:29:46"
$ws.Range("C$row").Value = "Monday"
$ws.Range("D$row").Value = "25"

# Drop the formatting we applied just to force text entry, so the new row
# does not pick up an extra cell style compared to the existing rows.
$textRange.ClearFormats()

# Numeric columns (city resale counts).
$ws.Range("E$row").Value = 122089
$ws.Range("F$row").Value = 133831
$ws.Range("G$row").Value = 162205
$ws.Range("H$row").Value = 133055
$ws.Range("I$row").Value = 177374
$ws.Range("J$row").Value = 114487
$ws.Range("K$row").Value = 201452
$ws.Range("L$row").Value = 225081
$ws.Range("M$row").Value = 175641
$ws.Range("N$row").Value = 103759
$ws.Range("O$row").Value = 39124
$ws.Range("P$row").Value = 33972
$ws.Range("Q$row").Value = 51834
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36198
$ws.Range("T$row").Value = -1
